# Applies the "updated revision data sheets" commit to the workbook.
# Only the "Skill Point Distributions" sheet has real content changes;
# all shared-string index shifts elsewhere happen automatically when the
# workbook is re-saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skill Point Distributions")

# Swords section (row 6 headers / row 7 values):
# "Overwhelm" perk was replaced by "Deadly Precision"
$ws.Range("M6").Value2 = "Deadly Precision"

# Row 7 point allocations changed for several Swords perks
$ws.Range("K7").Value2 = 1
$ws.Range("M7").Value2 = 2
$ws.Range("O7").Value2 = 1
$ws.Range("P7").Value2 = 1
$ws.Range("Q7").Value2 = 2

# Signs section (row 9 headers / row 10 values):
# Renamed Sign perks
$ws.Range("J9").Value2 = "The 5th Element"
$ws.Range("M9").Value2 = "Control Over The Power"
$ws.Range("Q9").Value2 = "Raw Power"

# Row 10 point allocations changed for several Sign perks
$ws.Range("E10").Value2 = 1
$ws.Range("F10").Value2 = 1
$ws.Range("I10").Value2 = 1
$ws.Range("Q10").Value2 = 1

# Update the active cell / selection on the sheet
$ws.Range("I15").Select()
